# Update workbook with newest numbers (as of 7th October 2024):
# appends 12 new records (rows 42-53, year 2024) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data: CVR, Year, TCV amount, Løsning, Opsagt dato (serial), Ny leverandør, Quarter, TCV_range
$cvr    = @("30510518","30972406","21825832","15516046","39525984","46481410","14036431","12759274","56577815","77144714","43268570","61552812")
$year   = @(2024,2024,2024,2024,2024,2024,2024,2024,2024,2024,2024,2024)
$amount = @(40925,59316,47491.62,50892,40956,59832,50400,55116,41460,40152,49752,53700)
$losn   = @("Visma Løn","BPO Visma Løn","Visma Løn","EasyCruit","Visma Time","Visma Løn","Visma Løn","Visma Løn","EasyCruit","Visma Løn","Visma Time","Visma Time")
$dato   = @(45301,45345,45370,45400,45329,45446,45447,45488,45516,45467,45526,45546)
$leverandor = @($null,$null,$null,$null,$null,$null,"Zenegy",$null,$null,$null,$null,$null)
$quarter = @("2024Q1","2024Q1","2024Q1","2024Q2","2024Q1","2024Q2","2024Q2","2024Q3","2024Q3","2024Q2","2024Q3","2024Q3")
$tcvrange = @("40000-60000","40000-60000","40000-60000","40000-60000","40000-60000","40000-60000","40000-60000","40000-60000","40000-60000","40000-60000","40000-60000","40000-60000")

$startRow = 42
$count = $cvr.Length

# The existing CVR column (A) stores numeric-looking values as TEXT (no leading
# apostrophe / quote-prefix style in the source file). A plain `.Value = "123"`
# assignment would be re-interpreted as a number, so instead we stage the text
# through a helper cell's TEXT() formula and paste-special the computed value,
# which preserves the shared-string (text) type without touching any styles.
$dateFmt = $ws.Cells.Item(2, 5).NumberFormat()

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $ws.Range("Z1").Formula = '=TEXT(' + $cvr[$i] + ',"0")'
    $ws.Range("Z1").Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4163)
}
$ws.Range("Z1").Clear()

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 4).Value = $losn[$i]
}

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    if ($leverandor[$i]) {
        $ws.Cells.Item($row, 7).Value = $leverandor[$i]
    }
}

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 8).Value = $quarter[$i]
}

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 9).Value = $tcvrange[$i]
}

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $year[$i]
    $ws.Cells.Item($row, 3).Value = $amount[$i]
    $ws.Cells.Item($row, 5).Value = $dato[$i]
    $ws.Cells.Item($row, 5).NumberFormat = $dateFmt
}
